$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "hhh not seen"
$ws.Cells.Item(10, 3).Value = "fixed"
$ws.Cells.Item(10, 4).Value = "`n            Passed"
$ws.Cells.Item(10, 5).Value = "2023-07-28 11:04:09"
$ws.Cells.Item(10, 6).Value = "2023-07-28 11:04:28"

# Match the bordered/bold/centered style used by column A on the prior rows
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
